$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.758234
$ws.Range("H2").Value = 8.274702
$ws.Range("I2").Value = 0.1921318935545868
$ws.Range("J2").Value = 0.1921318935545868
$ws.Range("M2").Value = 0.2992453333333333
$ws.Range("N2").Value = 0.897736
$ws.Range("O2").Value = 0.1430649508272797
$ws.Range("P2").Value = 0.1430649508272797
$ws.Range("Q2").Value = 0.8253886527413332
$ws.Range("R2").Value = 7.428497874672
$ws.Range("S2").Value = 0.0274873399037391
$ws.Range("T2").Value = 0.0274873399037391
$ws.Range("G3").Value = 2.758234
$ws.Range("H3").Value = 8.274702
$ws.Range("I3").Value = 0.1921318935545868
$ws.Range("J3").Value = 0.1921318935545868
$ws.Range("N3").Value = 3.4087
$ws.Range("O3").Value = 0.5432170458630915
$ws.Range("P3").Value = 0.5432170458630915
$ws.Range("Q3").Value = 3.133997411933333
$ws.Range("R3").Value = 28.20597670739999
$ws.Range("S3").Value = 0.1043693196328046
$ws.Range("T3").Value = 0.1043693196328046
$ws.Range("G4").Value = 2.758234
$ws.Range("H4").Value = 8.274702
$ws.Range("I4").Value = 0.1921318935545868
$ws.Range("J4").Value = 0.1921318935545868
$ws.Range("M4").Value = 0.656196
$ws.Range("N4").Value = 1.968588
$ws.Range("O4").Value = 0.3137180033096288
$ws.Range("P4").Value = 0.3137180033096288
$ws.Range("Q4").Value = 1.809942117864
$ws.Range("R4").Value = 16.289479060776
$ws.Range("S4").Value = 0.06027523401804311
$ws.Range("T4").Value = 0.06027523401804311
$ws.Range("G5").Value = 7.720664
$ws.Range("I5").Value = 0.5378027367579149
$ws.Range("J5").Value = 0.5378027367579149
$ws.Range("M5").Value = 0.2992453333333333
$ws.Range("N5").Value = 0.897736
$ws.Range("O5").Value = 0.1430649508272797
$ws.Range("P5").Value = 0.1430649508272797
$ws.Range("Q5").Value = 2.310372672234667
$ws.Range("R5").Value = 20.793354050112
$ws.Range("S5").Value = 0.07694072208904754
$ws.Range("T5").Value = 0.07694072208904754
$ws.Range("G6").Value = 7.720664
$ws.Range("I6").Value = 0.5378027367579149
$ws.Range("J6").Value = 0.5378027367579149
$ws.Range("N6").Value = 3.4087
$ws.Range("O6").Value = 0.5432170458630915
$ws.Range("P6").Value = 0.5432170458630915
$ws.Range("R6").Value = 78.95228213039999
$ws.Range("S6").Value = 0.2921436139187204
$ws.Range("T6").Value = 0.2921436139187204
$ws.Range("G7").Value = 7.720664
$ws.Range("I7").Value = 0.5378027367579149
$ws.Range("J7").Value = 0.5378027367579149
$ws.Range("M7").Value = 0.656196
$ws.Range("N7").Value = 1.968588
$ws.Range("O7").Value = 0.3137180033096288
$ws.Range("P7").Value = 0.3137180033096288
$ws.Range("Q7").Value = 5.066268834144
$ws.Range("R7").Value = 45.59641950729601
$ws.Range("S7").Value = 0.168718400750147
$ws.Range("T7").Value = 0.168718400750147
$ws.Range("G8").Value = 3.877042333333333
$ws.Range("H8").Value = 11.631127
$ws.Range("I8").Value = 0.2700653696874982
$ws.Range("J8").Value = 0.2700653696874982
$ws.Range("M8").Value = 0.2992453333333333
$ws.Range("N8").Value = 0.897736
$ws.Range("O8").Value = 0.1430649508272797
$ws.Range("P8").Value = 0.1430649508272797
$ws.Range("Q8").Value = 1.160186825385778
$ws.Range("R8").Value = 10.441681428472
$ws.Range("S8").Value = 0.03863688883449305
$ws.Range("T8").Value = 0.03863688883449305
$ws.Range("G9").Value = 3.877042333333333
$ws.Range("H9").Value = 11.631127
$ws.Range("I9").Value = 0.2700653696874982
$ws.Range("J9").Value = 0.2700653696874982
$ws.Range("N9").Value = 3.4087
$ws.Range("O9").Value = 0.5432170458630915
$ws.Range("P9").Value = 0.5432170458630915
$ws.Range("Q9").Value = 4.405224733877777
$ws.Range("R9").Value = 39.64702260489999
$ws.Range("S9").Value = 0.1467041123115665
$ws.Range("T9").Value = 0.1467041123115665
$ws.Range("G10").Value = 3.877042333333333
$ws.Range("H10").Value = 11.631127
$ws.Range("I10").Value = 0.2700653696874982
$ws.Range("J10").Value = 0.2700653696874982
$ws.Range("M10").Value = 0.656196
$ws.Range("N10").Value = 1.968588
$ws.Range("O10").Value = 0.3137180033096288
$ws.Range("P10").Value = 0.3137180033096288
$ws.Range("Q10").Value = 2.544099670964
$ws.Range("R10").Value = 22.896897038676
$ws.Range("S10").Value = 0.08472436854143869
$ws.Range("T10").Value = 0.08472436854143869
